# SCRUM board update for Sprint 2
# Re-triage the task list: everything that used to live in the "To-Do" (A)
# column (other than the still-open "aesthetically pleasing costs" task)
# is redistributed into "In Progress" (B) or "Done" (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Capture the current (pre-sprint) task text for each row in column A ---
$task3  = $ws.Cells.Item(3, 1).Value2   # Come up with the best method to store user information
$task5  = $ws.Cells.Item(5, 1).Value2   # Find a way to upload pictures to the website
$task6  = $ws.Cells.Item(6, 1).Value2   # Make the picture accessible to everyone who is splitting payment
$task7  = $ws.Cells.Item(7, 1).Value2   # Write a separate method to divide costs by percentage
$task8  = $ws.Cells.Item(8, 1).Value2   # Figure out how to send a text or email notifications over web apps
$task9  = $ws.Cells.Item(9, 1).Value2   # Write a function that divides costs
$task10 = $ws.Cells.Item(10, 1).Value2  # Integrate the Facebook JavaScript sdk
$task11 = $ws.Cells.Item(11, 1).Value2  # Test the website to find the best location for Facebook integration
$task12 = $ws.Cells.Item(12, 1).Value2  # Begin implementing whatever we come up with
$task13 = $ws.Cells.Item(13, 1).Value2  # Implement notification functionality

# --- Wipe out column A for the rows that move (row 4 stays put) ---
$ws.Range("A3").Clear()
$ws.Range("A5:A13").Clear()

# --- Drop each task into its new Sprint 2 column ---
$ws.Cells.Item(3, 3).Value  = $task3    # -> Done
$ws.Cells.Item(5, 3).Value  = $task5    # -> Done
$ws.Cells.Item(6, 3).Value  = $task6    # -> Done
$ws.Cells.Item(7, 3).Value  = $task7    # -> Done
$ws.Cells.Item(8, 2).Value  = $task8    # -> In Progress
$ws.Cells.Item(9, 2).Value  = $task9    # -> In Progress
$ws.Cells.Item(10, 2).Value = $task10   # -> In Progress
$ws.Cells.Item(11, 2).Value = $task11   # -> In Progress
$ws.Cells.Item(12, 3).Value = $task12   # -> Done
$ws.Cells.Item(13, 3).Value = $task13   # -> Done

# --- Row height tweaks for the rows whose wrapped text now needs more room ---
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 38.25

# --- Move the active selection to reflect where work left off ---
$ws.Range("C18").Select() | Out-Null
